$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.483.31'
$ws.Range("E2").Value = '  +1.48%  '

$ws.Range("D3").Value = '2.584.96'
$ws.Range("E3").Value = '  +0.25%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = "'507.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.51%  '

$ws.Range("D6").Value = "'153.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.52%  '

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = "'0.579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.75%  '

$ws.Range("D9").Value = '2.591.74'
$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("D10").Value = "'6.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.22%  '

$ws.Range("E11").Value = '  +1.00%  '

$ws.Range("D12").Value = "'0.348"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.95%  '

$ws.Range("E13").Value = '  +0.74%  '

$ws.Range("D14").Value = '3.040.78'
$ws.Range("E14").Value = '  +0.66%  '

$ws.Range("D15").Value = '60.475.50'
$ws.Range("E15").Value = '  +1.73%  '

$ws.Range("D16").Value = "'21.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.90%  '

$ws.Range("E17").Value = '  +2.31%  '

$ws.Range("D18").Value = '2.594.71'
$ws.Range("E18").Value = '  +0.63%  '

$ws.Range("D19").Value = "'4.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.38%  '

$ws.Range("D20").Value = "'345.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.67%  '

$ws.Range("D21").Value = "'10.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.70%  '

$ws.Range("D22").Value = "'6.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.56%  '

$ws.Range("D23").Value = "'0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.65%  '

$ws.Range("D24").Value = "'59.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.17%  '

$ws.Range("D25").Value = "'0.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.88%  '

$ws.Range("E26").Value = '  +1.07%  '

$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("D28").Value = '0.0₃0843'
$ws.Range("E28").Value = '  +0.56%  '

$ws.Range("D29").Value = "'7.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.16%  '

$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("D31").Value = "'19.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("D32").Value = "'153.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.23%  '

$ws.Range("E33").Value = '  -0.65%  '

$ws.Range("D34").Value = "'5.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.64%  '

$ws.Range("D35").Value = "'3.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.92%  '

$ws.Range("E36").Value = '  -0.46%  '

$ws.Range("E37").Value = '  +8.94%  '

$ws.Range("D38").Value = "'0.850"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.09%  '

$ws.Range("E39").Value = '  +2.61%  '

$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("D41").Value = "'35.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.26%  '

$ws.Range("D42").Value = "'295.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.62%  '

$ws.Range("E43").Value = '  -0.49%  '

$ws.Range("D44").Value = "'0.0995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.44%  '

$ws.Range("D45").Value = "'0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.22%  '

$ws.Range("E46").Value = '  -3.36%  '

$ws.Range("D47").Value = "'19.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.60%  '

$ws.Range("D48").Value = "'4.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.34%  '

$ws.Range("D49").Value = "'0.0233"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.26%  '

$ws.Range("E50").Value = '  +0.45%  '

$ws.Range("D51").Value = '2.000.82'
$ws.Range("E51").Value = '  +0.54%  '
